$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 5.040405773919509
$ws.Range("F3").Value = 5.904719908125438
$ws.Range("F4").Value = 6.508801240497247
$ws.Range("F5").Value = 32
$ws.Range("F6").Value = 7.245857245321472
$ws.Range("F7").Value = 5.809812808723561
$ws.Range("F9").Value = 0.8528893171779233
